# Discharge_Aug01.xlsx edit: "lots of discharge data"
#
# - stn3 sheet (sheet2.xml): append a third "new depth" data block
#   (rows 33-46), mirroring the existing "new velocity" block (rows 16-29).
# - stn1 sheet (sheet1.xml): selection moved; C37 un-shares its formula.
# - stn4 sheet (sheet3.xml): selection moved.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("stn1")
$ws3 = $wb.Worksheets.Item("stn3")
$ws4 = $wb.Worksheets.Item("stn4")

# ---------------------------------------------------------------------
# stn1 (sheet1.xml)
# ---------------------------------------------------------------------

# C37 becomes an explicit (non-shared) formula instead of continuing the
# C31:C38 shared-formula group.
$ws1.Range("C37").Formula = "=C23*2.54"

# ---------------------------------------------------------------------
# stn3 (sheet2.xml) - new "new depth" block, rows 33-46
# ---------------------------------------------------------------------

# Row 33: section header (bold), new shared string "new depth"
$ws3.Range("A33").Value = "new depth"
$ws3.Range("A33").Font.Bold = $true

# Row 34: column headers, reusing existing shared strings
$ws3.Range("A34").Value = "X"
$ws3.Range("B34").Value = "V"
$ws3.Range("C34").Value = "D"
$ws3.Range("D34").Value = "segment"
$ws3.Range("E34").Value = "Q"
$ws3.Range("F34").Value = "Qtotal"

# Column A (depth), rows 35-46
$aVals = @(0.7, 0.75, 0.8, 0.85, 0.9, 0.95, 1, 1.05, 1.1, 1.15, 1.2, 1.25)
for ($i = 0; $i -lt $aVals.Length; $i++) {
    $ws3.Cells.Item(35 + $i, 1).Value = $aVals[$i]
}

# Column B (velocity), rows 35-46
$bVals = @(0, 0, 0.14872000000000002, 0.1716, 0.0858, 0.08008, 0.0572, 0.08008, 0.10868, 0.0858, 0, 0)
for ($i = 0; $i -lt $bVals.Length; $i++) {
    $ws3.Cells.Item(35 + $i, 2).Value = $bVals[$i]
}

# Column C (new depth, cm) = old-block C * 2.54, shifted 17 rows
$ws3.Range("C35").Formula = "=C18*2.54"
$ws3.Range("C36:C46").Formula = "=C19*2.54"

# Column D (segment midpoint)
$ws3.Range("D35").Formula = "=A35"
$ws3.Range("D36").Formula = "=(A36+(A37-A36)/2)"
$ws3.Range("D37:D38").Formula = "=(A37+(A38-A37)/2)"
$ws3.Range("D39").Formula = "=(A39+(A40-A39)/2)"
$ws3.Range("D40:D46").Formula = "=(A40+(A41-A40)/2)"

# Column E (segment discharge)
$ws3.Range("E36").Formula = "=(D36-D35)*(B36)*C36"
$ws3.Range("E37").Formula = "=(D37-D36)*(B37)*C37"
$ws3.Range("E38:E46").Formula = "=(D38-D37)*(B38)*C38"

# Column F (total)
$ws3.Range("F35").Formula = "=SUM(E35:E53)"

# ---------------------------------------------------------------------
# Selections (saved per-sheet view state)
# ---------------------------------------------------------------------
[void]$ws1.Range("D36").Select()
[void]$ws3.Range("F35").Select()
[void]$ws4.Range("E40").Select()

# stn1 remains the active/visible tab, as in the original workbook.
[void]$ws1.Activate()
